$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    2  = -2
    10 = 0
    15 = 4
    18 = -1
    28 = 4
    32 = 0
    33 = -6
    34 = 2
    40 = -3
    42 = 1
    46 = 3
    47 = -2
    49 = -3
    51 = 4
    55 = 1
    60 = 1
    70 = 0
    71 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
